$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$origSheet = $wb.ActiveSheet

# Insert a new row above row 64, shifting existing rows (and the table)
# down by one.
$ws.Range("A64:K64").Insert(-4121)

# The newly inserted row picks up a blank default style; copy the
# formatting from the row directly below (the old row 64, now row 65) so
# it matches the rest of the leave-card table.
$ws.Range("A65:K65").Copy()
$ws.Range("A64:K64").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new leave-card entry: an undertime (UT) charge of 56
# minutes, recorded as an "Absence Undertime W/ Pay" debit.
$ws.Range("B64").Value = "UT(0-0-56)"
$ws.Range("D64").Value = 0.117
$ws.Range("G64").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# Grow the table to cover the newly inserted row.
$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A8:K136"))
$ws.Range("G136").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# Record the 56 minutes undertime on the CONVERTION sheet so its lookup
# table converts it into the equivalent-day value used above.
$conv = $wb.Worksheets.Item("CONVERTION")
$conv.Range("F3").Value = 56
$conv.Range("G3").Select()

# Restore the view to the working sheet / cell.
$origSheet.Activate()
$ws.Range("E73").Select()

$wb.RefreshAll()
$excel.CalculateFull()
